# TW edits for diagram and service name
#
# Re-positions/re-sizes shapes in the IBM MQ on EKS architecture diagram:
# the three "Public subnet" groups (box + icon), their "NAT gateway" labels,
# the "Boot node" label/icon, and the "EKS managed node group" background
# rectangle. Only a:off / a:ext values change; no text is edited.
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points while the
# canonical OOXML stores EMUs (914400 EMU = 72 pt). The point literals below
# were chosen so that, after the runtime's internal point<->EMU round trip,
# they reproduce the exact target EMU values from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)


# --- Left group ("Public subnet" / NAT gateway / Boot node), Availability Zone 1 ---

# Rectangle 67 (id 68) - "Public subnet" box (left group)
$shp = $s.Shapes.Item(5)
$shp.Left   = 162.72023782047242
$shp.Top    = 180.0
$shp.Width  = 175.54386146771654
$shp.Height = 151.2000046

# Graphic 68 (id 69) - public-subnet icon (left group)
$shp = $s.Shapes.Item(6)
$shp.Left   = 162.72023782047242
$shp.Top    = 180.0
$shp.Width  = 30.0
$shp.Height = 30.0

# TextBox 17 (id 74) - "NAT gateway" label (left group)
$shp = $s.Shapes.Item(7)
$shp.Left   = 207.29212598425198
$shp.Top    = 243.35787401574802
$shp.Width  = 86.4
$shp.Height = 20.599212598425197

# Graphic 35 (id 75) - NAT gateway icon (left group)
$shp = $s.Shapes.Item(8)
$shp.Left   = 232.08149606299213
$shp.Top    = 208.8
$shp.Width  = 36.0
$shp.Height = 36.0

# Graphic 60 (id 95) - boot-node icon (left group)
$shp = $s.Shapes.Item(11)
$shp.Left   = 232.08149606299213
$shp.Top    = 277.2
$shp.Width  = 36.0
$shp.Height = 36.0

# TextBox 16 (id 96) - "Boot node" label (left group)
$shp = $s.Shapes.Item(12)
$shp.Left   = 220.2407914015748
$shp.Top    = 309.6
$shp.Width  = 57.6000004
$shp.Height = 20.599212598425197


# --- Middle group ("Public subnet" / NAT gateway), Availability Zone 2 ---

# Rectangle 183 (id 184) - "Public subnet" box (middle group)
$shp = $s.Shapes.Item(16)
$shp.Left   = 480.0640157480315
$shp.Top    = 180.0
$shp.Width  = 175.6800004
$shp.Height = 151.2000046

# Graphic 184 (id 185) - public-subnet icon (middle group)
$shp = $s.Shapes.Item(17)
$shp.Left   = 480.063949607874
$shp.Top    = 180.0
$shp.Width  = 30.0
$shp.Height = 30.0

# TextBox 17 (id 186) - "NAT gateway" label (middle group)
$shp = $s.Shapes.Item(18)
$shp.Left   = 523.8412781425196
$shp.Top    = 243.36
$shp.Width  = 86.4
$shp.Height = 20.599212598425197

# Graphic 35 (id 187) - NAT gateway icon (middle group)
$shp = $s.Shapes.Item(19)
$shp.Left   = 549.7836220472441
$shp.Top    = 208.8
$shp.Width  = 36.0
$shp.Height = 36.0


# --- "EKS managed node group" background rectangle ---

# Rectangle 52 (id 53) - "EKS managed / node group" background rectangle
$shp = $s.Shapes.Item(26)
$shp.Left   = 135.61307526614172
$shp.Top    = 381.4036220472441
$shp.Width  = 792.0
$shp.Height = 230.4000016


# --- Right group ("Public subnet" / NAT gateway), Availability Zone 2 (second copy) ---

# Rectangle 13 (id 14) - "Public subnet" box (right group)
$shp = $s.Shapes.Item(34)
$shp.Left   = 725.4174803149606
$shp.Top    = 180.0
$shp.Width  = 175.6800004
$shp.Height = 151.2000046

# Graphic 14 (id 15) - public-subnet icon (right group)
$shp = $s.Shapes.Item(35)
$shp.Left   = 725.4174015748032
$shp.Top    = 180.0
$shp.Width  = 30.0
$shp.Height = 30.0

# TextBox 17 (id 16) - "NAT gateway" label (right group)
$shp = $s.Shapes.Item(36)
$shp.Left   = 768.894439748819
$shp.Top    = 243.35787401574802
$shp.Width  = 86.4
$shp.Height = 20.599212598425197

# Graphic 35 (id 17) - NAT gateway icon (right group)
$shp = $s.Shapes.Item(37)
$shp.Left   = 795.1371155141732
$shp.Top    = 208.8
$shp.Width  = 36.0
$shp.Height = 36.0
